# Clean-up of input tables:
#  - rename the worksheet from its temporary "updated" label back to the
#    default "Tabelle1" name
#  - move the live cell selection from T9 to O10 (last place the author
#    was working before saving)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Tabelle1"

$ws.Range("O10").Select()
